$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: phone number and quantity are stored as real numbers instead of text
$ws.Range("D4").Value = 508988751
$ws.Range("H4").Value = 4

# New row 5: a new order entry (phone number & quantity kept as text, like the
# original row 4 values were before they got normalized to numbers above)
$ws.Range("A5").Value = "ORD-20250301000448"
$ws.Range("B5").Value = "2025-03-01 00:04:48"
$ws.Range("C5").Value = "Steven"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0508988751"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "Jeans"
$ws.Range("F5").Value = "S"
$ws.Range("G5").Value = "Green"

$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "5"
$ws.Range("H5").Style = "Normal"

$ws.Range("I5").Value = "Dubai"
$ws.Range("J5").Value = "New"
